$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "WSL2 Ubuntu설치"
$ws.Range("E28").Value = "https://ropiens.tistory.com/155"

$ws.Range("D36").Value = "Deep semi-supervised learning (Basic and Algorithms)"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/337"

$ws.Range("D37").Value = "[Paper Review] HOTPOTQA: A Dataset for Diverse, Explainable  Multi-hop Question Answering"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1824&mod=document&pageid=1"
